{"js": "const replacements = [\n  [\"2024-03-28 Thursday\", \"2024-03-29 Friday\"],\n  [\"540\u00f73=\", \"148\u00f73=\"],\n  [\"686\u00f72=\", \"640\u00f73=\"],\n  [\"703\u00f75=\", \"281\u00f73=\"],\n  [\"639\u00f79=\", \"670\u00f73=\"],\n  [\"219\u00f72=\", \"435\u00f72=\"],\n  [\"219\u00f78=\", \"862\u00f72=\"],\n  [\"706\u00f74=\", \"167\u00f72=\"],\n  [\"582\u00f79=\", \"685\u00f77=\"],\n  [\"779\u00f74=\", \"944\u00f79=\"],\n  [\"338\u00f74=\", \"692\u00f78=\"],\n  [\"755\u00f78=\", \"606\u00f76=\"],\n  [\"419\u00f78=\", \"836\u00f72=\"],\n  [\"509\u00f73=\", \"442\u00f79=\"],\n  [\"186\u00f72=\", \"833\u00f74=\"],\n  [\"229\u00f77=\", \"257\u00f72=\"],\n  [\"861\u00f73=\", \"489\u00f75=\"],\n  [\"464\u00f72=\", \"737\u00f78=\"],\n  [\"540\u00f72=\", \"771\u00f72=\"],\n  [\"939\u00f75=\", \"433\u00f76=\"],\n  [\"189\u00f75=\", \"994\u00f76=\"],\n  [\"442\u00f74=\", \"838\u00f73=\"],\n  [\"250\u00f73=\", \"485\u00f73=\"],\n  [\"273\u00f77=\", \"225\u00f75=\"],\n  [\"198\u00f77=\", \"908\u00f74=\"],\n  [\"300\u00f77=\", \"165\u00f79=\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    ,@(\"2024-03-28 Thursday\", \"2024-03-29 Friday\")\n    ,@(\"540\u00f73=\", \"148\u00f73=\")\n    ,@(\"686\u00f72=\", \"640\u00f73=\")\n    ,@(\"703\u00f75=\", \"281\u00f73=\")\n    ,@(\"639\u00f79=\", \"670\u00f73=\")\n    ,@(\"219\u00f72=\", \"435\u00f72=\")\n    ,@(\"219\u00f78=\", \"862\u00f72=\")\n    ,@(\"706\u00f74=\", \"167\u00f72=\")\n    ,@(\"582\u00f79=\", \"685\u00f77=\")\n    ,@(\"779\u00f74=\", \"944\u00f79=\")\n    ,@(\"338\u00f74=\", \"692\u00f78=\")\n    ,@(\"755\u00f78=\", \"606\u00f76=\")\n    ,@(\"419\u00f78=\", \"836\u00f72=\")\n    ,@(\"509\u00f73=\", \"442\u00f79=\")\n    ,@(\"186\u00f72=\", \"833\u00f74=\")\n    ,@(\"229\u00f77=\", \"257\u00f72=\")\n    ,@(\"861\u00f73=\", \"489\u00f75=\")\n    ,@(\"464\u00f72=\", \"737\u00f78=\")\n    ,@(\"540\u00f72=\", \"771\u00f72=\")\n    ,@(\"939\u00f75=\", \"433\u00f76=\")\n    ,@(\"189\u00f75=\", \"994\u00f76=\")\n    ,@(\"442\u00f74=\", \"838\u00f73=\")\n    ,@(\"250\u00f73=\", \"485\u00f73=\")\n    ,@(\"273\u00f77=\", \"225\u00f75=\")\n    ,@(\"198\u00f77=\", \"908\u00f74=\")\n    ,@(\"300\u00f77=\", \"165\u00f79=\")\n)\n\nforeach ($pair in $replacements) {\n    $findText = $pair[0]\n    $replaceText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n"}
